$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: "time_taken" header + per-row timestamp metadata
$ws.Range("F1").Value = "time_taken"

# Mirror the header style (bold, centered, bordered) used by B1:E1 onto F1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data cells F2:F31 - text timestamp values (no explicit style, like B2:E31)
$ws.Range("F2").Value = "2021-10-05 13:40:34.160376"
$ws.Range("F3").Value = "2021-10-05 13:40:34.160386"
$ws.Range("F4").Value = "2021-10-05 13:40:34.160389"
$ws.Range("F5").Value = "2021-10-05 13:40:34.160392"
$ws.Range("F6").Value = "2021-10-05 13:40:34.160395"
$ws.Range("F7").Value = "2021-10-05 13:40:34.160398"
$ws.Range("F8").Value = "2021-10-05 13:40:34.160400"
$ws.Range("F9").Value = "2021-10-05 13:40:34.160403"
$ws.Range("F10").Value = "2021-10-05 13:40:34.160406"
$ws.Range("F11").Value = "2021-10-05 13:40:34.160408"
$ws.Range("F12").Value = "2021-10-05 13:40:34.160411"
$ws.Range("F13").Value = "2021-10-05 13:40:34.160413"
$ws.Range("F14").Value = "2021-10-05 13:40:34.160416"
$ws.Range("F15").Value = "2021-10-05 13:40:34.160418"
$ws.Range("F16").Value = "2021-10-05 13:40:34.160421"
$ws.Range("F17").Value = "2021-10-05 13:40:34.160423"
$ws.Range("F18").Value = "2021-10-05 13:40:34.160426"
$ws.Range("F19").Value = "2021-10-05 13:40:34.160429"
$ws.Range("F20").Value = "2021-10-05 13:40:34.160432"
$ws.Range("F21").Value = "2021-10-05 13:40:34.160435"
$ws.Range("F22").Value = "2021-10-05 13:40:34.160437"
$ws.Range("F23").Value = "2021-10-05 13:40:34.160440"
$ws.Range("F24").Value = "2021-10-05 13:40:34.160442"
$ws.Range("F25").Value = "2021-10-05 13:40:34.160445"
$ws.Range("F26").Value = "2021-10-05 13:40:34.160448"
$ws.Range("F27").Value = "2021-10-05 13:40:34.160451"
$ws.Range("F28").Value = "2021-10-05 13:40:34.160453"
$ws.Range("F29").Value = "2021-10-05 13:40:34.160456"
$ws.Range("F30").Value = "2021-10-05 13:40:34.160458"
$ws.Range("F31").Value = "2021-10-05 13:40:34.160461"
